$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = 971.2554529537363
$ws.Range("C23").Value = 75251889187.105
$ws.Range("E23").Value = 672.81932515247
$ws.Range("E60").Value = 607.0440914077899
$ws.Range("E276").Value = 554.9665658524307
$ws.Range("C326").Value = 27568366764.73001
$ws.Range("E326").Value = 550.8860153855919
$ws.Range("E366").Value = 968.0084813712562
$ws.Range("F366").Value = 2134.071498031071
$ws.Range("C407").Value = 38405119201.1
$ws.Range("E407").Value = 759.5055166232634
$ws.Range("F407").Value = 1674.405861947647
$ws.Range("C418").Value = 86092205729.145
$ws.Range("E418").Value = 949.9260775150279
$ws.Range("F418").Value = 2094.207030489631
$ws.Range("E474").Value = 585.0274667584168
$ws.Range("C815").Value = 66479177033.20326
$ws.Range("E815").Value = 361.9885619607754
$ws.Range("F815").Value = 798.0399836987255
$ws.Range("C863").Value = 30700167574.87013
$ws.Range("C886").Value = 75497054015.65343
$ws.Range("E886").Value = 597.581948306792
$ws.Range("F886").Value = 1317.429163237154
$ws.Range("C894").Value = 211867781160.6168
$ws.Range("E894").Value = 443.1749075544793
$ws.Range("F894").Value = 977.0234011946051
$ws.Range("C907").Value = 33804157477.81094
$ws.Range("E907").Value = 593.3555059633622
$ws.Range("F907").Value = 1308.111548446828
$ws.Range("C936").Value = 66724266888.81259
$ws.Range("E936").Value = 555.1321900075351
$ws.Range("F936").Value = 1223.844426090612
$ws.Range("C944").Value = 198980412292.8443
$ws.Range("E944").Value = 410.9450461882761
$ws.Range("F944").Value = 905.9694488266734
$ws.Range("C953").Value = 44472448824.8619
$ws.Range("E953").Value = 327.3259050133436
$ws.Range("F953").Value = 721.6226901924174
$ws.Range("C954").Value = 23299832052.79658
$ws.Range("E954").Value = 439.3975419495976
$ws.Range("F954").Value = 968.695820982083
$ws.Range("C955").Value = 40132962486.93383
$ws.Range("E955").Value = 354.3344843607576
$ws.Range("F955").Value = 781.1658042217263
$ws.Range("C956").Value = 38579055002.29417
$ws.Range("E956").Value = 180.420811003219
$ws.Range("F956").Value = 397.7557199376967
$ws.Range("C957").Value = 29668002872.75095
$ws.Range("E957").Value = 540.4069586914067
$ws.Range("F957").Value = 1191.381181131075
$ws.Range("C958").Value = 9862371713.702158
$ws.Range("E958").Value = 236.1576859989763
$ws.Range("F958").Value = 520.6332345533432
$ws.Range("C959").Value = 1692476585.791443
$ws.Range("E959").Value = 329.3042344408528
$ws.Range("F959").Value = 725.9841152483042
$ws.Range("C960").Value = 95205331494.32219
$ws.Range("E960").Value = 379.7418250053872
$ws.Range("F960").Value = 837.1788274068766
$ws.Range("C961").Value = 38841786515.80583
$ws.Range("E961").Value = 325.4168193115178
$ws.Range("F961").Value = 717.4139198541723
$ws.Range("C963").Value = 16725866015.13769
$ws.Range("E963").Value = 281.118885400509
$ws.Range("F963").Value = 619.7546947539621
$ws.Range("C964").Value = 1601968448.482061
$ws.Range("E964").Value = 82.62625473885981
$ws.Range("F964").Value = 182.1578411972904
$ws.Range("C965").Value = 43209972681.16507
$ws.Range("E965").Value = 248.4140419805541
$ws.Range("F965").Value = 547.6535969503295
$ws.Range("C966").Value = 62955464856.34379
$ws.Range("E966").Value = 700.2293420225965
$ws.Range("F966").Value = 1543.725607423016
$ws.Range("C967").Value = 19858526191.73903
$ws.Range("E967").Value = 365.1310717101807
$ws.Range("F967").Value = 804.9679606922645
$ws.Range("C968").Value = 48559338876.30147
$ws.Range("E968").Value = 765.0125517299552
$ws.Range("F968").Value = 1686.546671543859
$ws.Range("C969").Value = 35198014991.9602
$ws.Range("E969").Value = 344.2504722315819
$ws.Range("F969").Value = 758.9345910817453
$ws.Range("C970").Value = 6991457674.526141
$ws.Range("E970").Value = 339.071280759815
$ws.Range("F970").Value = 747.5165455630881
$ws.Range("C971").Value = 10536438538.04916
$ws.Range("E971").Value = 284.212694023066
$ws.Range("F971").Value = 626.5753052432514
$ws.Range("C972").Value = 981543254.49993
$ws.Range("E972").Value = 93.56398573966675
$ws.Range("F972").Value = 206.2711629616693
$ws.Range("C973").Value = 43986265264.57739
$ws.Range("E973").Value = 418.7958358864679
$ws.Range("F973").Value = 923.2772997953072
$ws.Range("C974").Value = 19764288224.18655
$ws.Range("E974").Value = 348.7792758604293
$ws.Range("F974").Value = 768.9187915619024
$ws.Range("C975").Value = 52330770771.55927
$ws.Range("E975").Value = 709.1750035757959
$ws.Range("F975").Value = 1563.4472128832
$ws.Range("C976").Value = 26512097215.38089
$ws.Range("E976").Value = 403.0333162972099
$ws.Range("F976").Value = 888.5272491088289
$ws.Range("C977").Value = 9723478626.075483
$ws.Range("E977").Value = 410.8932782531201
$ws.Range("F977").Value = 905.8553212368286
$ws.Range("C978").Value = 36767910527.61337
$ws.Range("E978").Value = 295.7694987253039
$ws.Range("F978").Value = 652.0534368898051
$ws.Range("C979").Value = 26830851576.61391
$ws.Range("E979").Value = 627.1090494855044
$ws.Range("F979").Value = 1382.524610495743
$ws.Range("C980").Value = 20451669834.83244
$ws.Range("E980").Value = 556.526385664284
$ws.Range("F980").Value = 1226.918069835481
$ws.Range("C981").Value = 1673949742.494239
$ws.Range("E981").Value = 99.54727030907524
$ws.Range("F981").Value = 219.4619121233873
$ws.Range("C982").Value = 13881727568.19452
$ws.Range("E982").Value = 216.8985671029392
$ws.Range("F982").Value = 478.1745810351399
$ws.Range("C983").Value = 19779568388.30058
$ws.Range("E983").Value = 572.7033453903583
$ws.Range("F983").Value = 1262.581795247584
$ws.Range("C984").Value = 13148226826.05747
$ws.Range("E984").Value = 316.8971998522767
$ws.Range("F984").Value = 698.6315667943292
$ws.Range("C985").Value = 25125058770.06453
$ws.Range("E985").Value = 187.1067507172949
$ws.Range("F985").Value = 412.4955426313484
$ws.Range("C986").Value = 68056796722.34452
$ws.Range("E986").Value = 561.5751834254621
$ws.Range("F986").Value = 1238.048649379774
$ws.Range("C987").Value = 26517188419.56076
$ws.Range("E987").Value = 317.1369182107484
$ws.Range("F987").Value = 699.1600498874159
$ws.Range("C988").Value = 9767671035.329494
$ws.Range("E988").Value = 150.0392427173927
$ws.Range("F988").Value = 330.7765144947639
$ws.Range("C989").Value = 72598815899.25462
$ws.Range("E989").Value = 313.7290400661198
$ws.Range("F989").Value = 691.6470417297677
$ws.Range("C990").Value = 3197524198.504057
$ws.Range("E990").Value = 384.7522429218327
$ws.Range("F990").Value = 848.2247947454724
$ws.Range("C991").Value = 22632517255.58004
$ws.Range("E991").Value = 229.6738295771612
$ws.Range("F991").Value = 506.3389246858096
$ws.Range("C992").Value = 2112685384.998585
$ws.Range("E992").Value = 124.1708539222332
$ws.Range("F992").Value = 273.7470645569554
$ws.Range("C993").Value = 20934719722.75175
$ws.Range("E993").Value = 269.8899187842088
$ws.Range("F993").Value = 594.9993149516666
$ws.Range("C994").Value = 183315608049.2028
$ws.Range("E994").Value = 384.2331522582539
$ws.Range("F994").Value = 847.0804074685467
$ws.Range("C995").Value = 25774258851.83581
$ws.Range("E995").Value = 684.4785759384282
$ws.Range("F995").Value = 1509.001468513859
$ws.Range("C996").Value = 29776091745.55735
$ws.Range("E996").Value = 290.5581391622026
$ws.Range("F996").Value = 640.5644735969919
$ws.Range("C997").Value = 1451738.276800026
$ws.Range("E997").Value = 0.5736012966832797
$ws.Range("F997").Value = 1.264561418667959
$ws.Range("C998").Value = 11482373688.1693
$ws.Range("E998").Value = 100.3331016028039
$ws.Range("F998").Value = 221.1943557935415
$ws.Range("C999").Value = 32404283369.72654
$ws.Range("E999").Value = 530.8447668703154
$ws.Range("F999").Value = 1170.300373042297
$ws.Range("C1000").Value = 48897052155.93681
$ws.Range("E1000").Value = 861.2665877544873
$ws.Range("F1000").Value = 1898.748319363543
$ws.Range("C1001").Value = 37169259410.77391
$ws.Range("E1001").Value = 892.1304380091231
$ws.Range("F1001").Value = 1966.790763634913
